$d = $word.ActiveDocument

# --- Edit 1: "Ingeniero en sistemas desarrollador" (one run)
#             -> "Técnico" + " en sistemas desarrollador" (two runs)
# Word's own save-time run consolidation would merge two adjacent,
# identically-formatted runs back into one if we just set text in place.
# Briefly tracking the change and then accepting each recorded revision
# individually (rather than Document.Revisions.AcceptAll, which rewrites
# the whole run tree) keeps the insertion as a separate run while leaving
# the rest of the document untouched.
$d.TrackRevisions = $true

$rng = $d.Content
$found1 = $rng.Find.Execute("Ingeniero en sistemas desarrollador", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)
if ($found1) {
    $rng.Text = "Técnico"
    $rng2 = $d.Range($rng.End, $rng.End)
    $rng2.InsertAfter(" en sistemas desarrollador")
}

$d.TrackRevisions = $false
while ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}

# --- Edit 2: "Ingeniero mecatr" + "ónico" (two runs, split only by a
#             language-tag rPr) -> "Técnico en electrónica" (one run)
$rng3 = $d.Content
$found2 = $rng3.Find.Execute("Ingeniero mecatrónico", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
if ($found2) {
    $rng3.Text = "Técnico en electrónica"
}
